# Apply updated average-price / profit figures per the scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 557.5714
$ws.Range("I29").Value = 213.75
$ws.Range("J29").Value = 1016
$ws.Range("K29").Value = 641.25
$ws.Range("L29").Value = 3048
$ws.Range("M29").Value = -360.25
$ws.Range("N29").Value = -3610

$ws.Range("H33").Value = 776.85
$ws.Range("I33").Value = 846.625
$ws.Range("J33").Value = 497.75
$ws.Range("K33").Value = 846.625
$ws.Range("L33").Value = 497.75
$ws.Range("M33").Value = -617.625
$ws.Range("N33").Value = -955.75

$ws.Range("H58").Value = 5261.375
$ws.Range("J58").Value = 7599.8
$ws.Range("L58").Value = 22799.4
$ws.Range("N58").Value = -23099.4

$ws.Range("H70").Value = 29168094
$ws.Range("I70").Value = 41667650
$ws.Range("J70").Value = 23811142
$ws.Range("K70").Value = 125002950
$ws.Range("L70").Value = 71433426
$ws.Range("M70").Value = -125002680
$ws.Range("N70").Value = -71433966

$ws.Range("H73").Value = 29168094
$ws.Range("I73").Value = 41667650
$ws.Range("J73").Value = 23811142
$ws.Range("K73").Value = 125002950
$ws.Range("L73").Value = 71433426
$ws.Range("M73").Value = -125002014
$ws.Range("N73").Value = -71435298

$ws.Range("H86").Value = 77162340
$ws.Range("I86").Value = 111112960
$ws.Range("K86").Value = 111112960
$ws.Range("M86").Value = -111111837

$ws.Range("H87").Value = 59993
$ws.Range("J87").Value = 59993
$ws.Range("L87").Value = 59993
$ws.Range("N87").Value = -62489

$ws.Range("H89").Value = 77162340
$ws.Range("I89").Value = 111112960
$ws.Range("K89").Value = 555564800
$ws.Range("M89").Value = -555559184

$ws.Range("H90").Value = 59993
$ws.Range("J90").Value = 59993
$ws.Range("L90").Value = 179979
$ws.Range("N90").Value = -192459

$ws.Range("H141").Value = 5873.75
$ws.Range("I141").Value = 2831.6667
$ws.Range("K141").Value = 8495.000100000001
$ws.Range("M141").Value = -3315.000100000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H111").Value = 48512.332
$ws.Range("J111").Value = 48512.332
$ws.Range("L111").Value = 48512.332

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 66670588
$ws.Range("I86").Value = 3006.5715
$ws.Range("J86").Value = 125004730
$ws.Range("K86").Value = 3006.5715
$ws.Range("L86").Value = 125004730
$ws.Range("M86").Value = -1883.5715
$ws.Range("N86").Value = -125006976

$ws.Range("H89").Value = 66670588
$ws.Range("I89").Value = 3006.5715
$ws.Range("J89").Value = 125004730
$ws.Range("K89").Value = 15032.8575
$ws.Range("L89").Value = 625023650
$ws.Range("M89").Value = -9416.8575
$ws.Range("N89").Value = -625034882

$ws.Range("H105").Value = 4812.6665
$ws.Range("I105").Value = 2561.8
$ws.Range("K105").Value = 2561.8
$ws.Range("M105").Value = -814.8000000000002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7400.8716
$ws.Range("I31").Value = 3294.7307
$ws.Range("J31").Value = 15613.154
$ws.Range("K31").Value = 3294.7307
$ws.Range("L31").Value = 15613.154
$ws.Range("M31").Value = -2999.7307
$ws.Range("N31").Value = -16203.154

$ws.Range("H34").Value = 7400.8716
$ws.Range("I34").Value = 3294.7307
$ws.Range("J34").Value = 15613.154
$ws.Range("K34").Value = 3294.7307
$ws.Range("L34").Value = 15613.154
$ws.Range("M34").Value = -3092.7307
$ws.Range("N34").Value = -16017.154

$ws.Range("H58").Value = 12506626
$ws.Range("I58").Value = 55558164
$ws.Range("K58").Value = 55558164
$ws.Range("M58").Value = -55557961

$ws.Range("H62").Value = 25004032
$ws.Range("I62").Value = 31252540
$ws.Range("K62").Value = 31252540
$ws.Range("M62").Value = -31251916

$ws.Range("H65").Value = 25004032
$ws.Range("I65").Value = 31252540
$ws.Range("K65").Value = 156262700
$ws.Range("M65").Value = -156259580

$ws.Range("H86").Value = 7816485
$ws.Range("I86").Value = 10420332
$ws.Range("K86").Value = 10420332
$ws.Range("M86").Value = -10419209

$ws.Range("H89").Value = 7816485
$ws.Range("I89").Value = 10420332
$ws.Range("K89").Value = 52101660
$ws.Range("M89").Value = -52096044

$ws.Range("H122").Value = 4278.375
$ws.Range("I122").Value = 4262.6
$ws.Range("J122").Value = 4304.6665
$ws.Range("K122").Value = 12787.8
$ws.Range("L122").Value = 12913.9995
$ws.Range("M122").Value = -10337.8
$ws.Range("N122").Value = -17813.9995

$ws.Range("H136").Value = 12506626
$ws.Range("I136").Value = 55558164
$ws.Range("K136").Value = 166674492
$ws.Range("M136").Value = -166671942

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 44.363636
$ws.Range("I38").Value = 40.8
$ws.Range("K38").Value = 122.4
$ws.Range("M38").Value = 224.6

$ws.Range("H107").Value = 18182522
$ws.Range("J107").Value = 20000726
$ws.Range("L107").Value = 60002178
$ws.Range("N107").Value = -60006018

$ws.Range("H113").Value = 3025.48
$ws.Range("I113").Value = 851.375
$ws.Range("J113").Value = 4048.5881
$ws.Range("K113").Value = 2554.125
$ws.Range("L113").Value = 12145.7643
$ws.Range("M113").Value = -384.125
$ws.Range("N113").Value = -16485.7643

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2857.6428
$ws.Range("I80").Value = 2792.6667
$ws.Range("K80").Value = 2792.6667
$ws.Range("M80").Value = -1794.6667

$ws.Range("H83").Value = 2857.6428
$ws.Range("I83").Value = 2792.6667
$ws.Range("K83").Value = 13963.3335
$ws.Range("M83").Value = -8971.333500000001

$ws.Range("H122").Value = 32113.639
$ws.Range("I122").Value = 39521.145
$ws.Range("K122").Value = 118563.435
$ws.Range("M122").Value = -116113.435

$ws.Range("H126").Value = 2269.0667
$ws.Range("I126").Value = 2149.2856
$ws.Range("K126").Value = 6447.8568
$ws.Range("M126").Value = -3977.8568

$ws.Range("H132").Value = 6237.353
$ws.Range("I132").Value = 3067.5715
$ws.Range("J132").Value = 8456.200000000001
$ws.Range("K132").Value = 9202.7145
$ws.Range("L132").Value = 25368.6
$ws.Range("M132").Value = -6672.7145
$ws.Range("N132").Value = -30428.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3737.2917
$ws.Range("I7").Value = 2649.7222
$ws.Range("J7").Value = 7000
$ws.Range("K7").Value = 2649.7222
$ws.Range("L7").Value = 7000
$ws.Range("M7").Value = -2537.7222
$ws.Range("N7").Value = -7224

$ws.Range("H46").Value = 1780.7693
$ws.Range("I46").Value = 685.75
$ws.Range("J46").Value = 2267.4443
$ws.Range("K46").Value = 685.75
$ws.Range("L46").Value = 2267.4443
$ws.Range("M46").Value = -497.75
$ws.Range("N46").Value = -2643.4443

$ws.Range("H126").Value = 3737.2917
$ws.Range("I126").Value = 2649.7222
$ws.Range("J126").Value = 7000
$ws.Range("K126").Value = 7949.1666
$ws.Range("L126").Value = 21000
$ws.Range("M126").Value = -5479.1666
$ws.Range("N126").Value = -25940

$ws.Range("H136").Value = 9331.771000000001
$ws.Range("I136").Value = 4298.5454
$ws.Range("J136").Value = 11638.667
$ws.Range("K136").Value = 12895.6362
$ws.Range("L136").Value = 34916.001
$ws.Range("M136").Value = -10345.6362
$ws.Range("N136").Value = -40016.001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7692.222
$ws.Range("I62").Value = 5946.8
$ws.Range("K62").Value = 5946.8
$ws.Range("M62").Value = -5322.8

$ws.Range("H65").Value = 7692.222
$ws.Range("I65").Value = 5946.8
$ws.Range("K65").Value = 29734
$ws.Range("M65").Value = -26614

$ws.Range("H100").Value = 430.7143
$ws.Range("I100").Value = 430.7143
$ws.Range("K100").Value = 861.4286
$ws.Range("M100").Value = -320.4286

$ws.Range("H113").Value = 1096.7715
$ws.Range("I113").Value = 754
$ws.Range("K113").Value = 2262
$ws.Range("M113").Value = -92

# ARM row 111 previously had no LeveProfitHQ (N) value; the refresh populates it.
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("N111").Value = -56692.332

Write-Output "Applied $([string]200) cell updates across 8 sheets."
